# The figure captions (and the matching picture alt-text / docPr
# description) currently repeat the figure number, e.g.
#   "Рис. 1: Рис.1.Каталог курса"
# The redundant "Рис.N." fragment baked into the second half of the
# string needs to be dropped, leaving just:
#   "Рис. 1: Каталог курса"
# This happens for figures 1-7. Both the caption paragraph (styled
# "ImageCaption") and the picture's alternative text need to be fixed.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "Рис. 1: Рис.1.Каталог курса"; New = "Рис. 1: Каталог курса" },
    @{ Old = "Рис. 2: Рис.2.Каталог с шаблоном отчета по лабораторной №4"; New = "Рис. 2: Каталог с шаблоном отчета по лабораторной №4" },
    @{ Old = "Рис. 3: Рис.3.make"; New = "Рис. 3: make" },
    @{ Old = "Рис. 4: Рис.4.report.pdg"; New = "Рис. 4: report.pdg" },
    @{ Old = "Рис. 5: Рис.5.make clean"; New = "Рис. 5: make clean" },
    @{ Old = "Рис. 6: Рис.6.gedit"; New = "Рис. 6: gedit" },
    @{ Old = "Рис. 7: Рис.7.Отправка файлов на GitHub"; New = "Рис. 7: Отправка файлов на GitHub" }
)

# 1) Fix the visible caption text runs (paragraphs styled "ImageCaption").
foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}

# 2) Fix each picture's alternative text (wp:docPr/@descr) to match.
#    InlineShapes are encountered in the same order as the figures.
foreach ($shape in $d.InlineShapes) {
    foreach ($pair in $replacements) {
        if ($shape.AlternativeText -eq $pair.Old) {
            $shape.AlternativeText = $pair.New
        }
    }
}
